$d = $word.ActiveDocument

# Locate the "Characters" heading paragraph (first paragraph in this
# document, a centered, bold + underlined title) and insert a new,
# separate run reading "List of " immediately before the existing
# "Characters" run, carrying the same bold/bCs/underline formatting, so
# the heading reads "List of Characters".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Characters") {
        $target = $p
        break
    }
}

$r = $target.Range
$r.Collapse(1)  # wdCollapseStart

# Insert raw OOXML so the new text becomes its own <w:r> (matching the
# target run-for-run) instead of being folded into the neighbouring run
# that merely happens to share the same formatting.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">List of </w:t></w:r></w:p>'
$r.InsertXML($xml)
